$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A43").Value = "Alessio Farinati"
$ws.Range("B43").Value = "Daniele  Dalbosco | iMontagna"
$ws.Range("C43").Value = "Andrea Conzatti | FC Savignano"
$ws.Range("D43").Value = "ALESSIO FARINATI | Pinguini Trentini"
$ws.Range("E43").Value = "MARCO HEIDEMPERGHER | U.S. Guarna"
$ws.Range("F43").Value = "Federico Rippa | Vigili del Fusto"
